$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the third column of data (column F, input size 10000) for each trial row
$ws.Range("F6").Value = 78
$ws.Range("F7").Value = 76
$ws.Range("F8").Value = 33
$ws.Range("F9").Value = 34
$ws.Range("F10").Value = 33
$ws.Range("F11").Value = 10
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 178
$ws.Range("F17").Value = 186
$ws.Range("F18").Value = 125
$ws.Range("F19").Value = 173
$ws.Range("F20").Value = 168
$ws.Range("F21").Value = 10
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("F24").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 33
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 21
$ws.Range("F29").Value = 20
$ws.Range("F30").Value = 22
$ws.Range("F31").Value = 11
$ws.Range("F32").Value = 3
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = 2
$ws.Range("F35").Value = 2
$ws.Range("F36").Value = 13
$ws.Range("F37").Value = 11
$ws.Range("F38").Value = 2
$ws.Range("F39").Value = 1
$ws.Range("F40").Value = 1
$ws.Range("F41").Value = 144
$ws.Range("F42").Value = 3
$ws.Range("F43").Value = 3
$ws.Range("F44").Value = 3
$ws.Range("F45").Value = 3
$ws.Range("F46").Value = 7
$ws.Range("F47").Value = 2
$ws.Range("F48").Value = 2
$ws.Range("F49").Value = 1
$ws.Range("F50").Value = 1
$ws.Range("F51").Value = 3
$ws.Range("F52").Value = 2
$ws.Range("F53").Value = 1
$ws.Range("F54").Value = 2
$ws.Range("F55").Value = 2
$ws.Range("F56").Value = 3
$ws.Range("F57").Value = 1
$ws.Range("F58").Value = 4
$ws.Range("F59").Value = 1
$ws.Range("F60").Value = 0
$ws.Range("F61").Value = 1
$ws.Range("F62").Value = 1
$ws.Range("F63").Value = 1
$ws.Range("F64").Value = 1
$ws.Range("F65").Value = 1
$ws.Range("F66").Value = 11
$ws.Range("F67").Value = 2
$ws.Range("F68").Value = 1
$ws.Range("F69").Value = 2
$ws.Range("F70").Value = 2
$ws.Range("F71").Value = 6
$ws.Range("F72").Value = 1
$ws.Range("F73").Value = 1
$ws.Range("F74").Value = 1
$ws.Range("F75").Value = 2
$ws.Range("F76").Value = 5
$ws.Range("F77").Value = 2
$ws.Range("F78").Value = 2
$ws.Range("F79").Value = 8
$ws.Range("F80").Value = 10
$ws.Range("F81").Value = 3
$ws.Range("F82").Value = 2
$ws.Range("F83").Value = 1
$ws.Range("F84").Value = 1
$ws.Range("F85").Value = 1

# Update the active sheet view / selection to match the saved state
$ws.Application.ActiveWindow.ScrollRow = 68
$ws.Range("F86").Select()
